$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Rows 8 and 9 swap which product line they describe (description, barcode,
# category and brand all move to the other row). The Pantene line keeps its
# original numbers (G/H/I/J/L), while the Bazaar line - now on row 9 - records
# one extra sale: SalesQuantity 2 -> 3 and Turnover 1.86 -> 2.79. The
# grand-total row (35) is updated to match the new column sums
# (K: 351 -> 352, L: 443.11 -> 444.04).
# ---------------------------------------------------------------------------

# --- Row 8 becomes the "Pantene Shampoo Repair & Protect 360ml" line -------
$ws.Range("E8").Value = "Pantene® Shampoo Repair & Protect 360ml"

$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "4015600948016"
$ws.Range("F8").Font.Name = "Avenir Next"
$ws.Range("F8").HorizontalAlignment = -4131

$ws.Range("A8").Value = "Πωλήσεις Έκπτωση 1"
$ws.Range("G8").Value = 3.98
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 30
$ws.Range("J8").Value = "Pantene"
$ws.Range("L8").Value = 4.39

# --- Row 9 becomes the "Bazaar Ygro Patomatos Ultra Anthi Paschalias" line --
$ws.Range("E9").Value = "Bazaar® Υγρό Πατώματος Ultra Άνθη Πασχαλιάς 1ltr"

$ws.Range("F9").NumberFormat = "@"
$ws.Range("F9").Value = "5208086416820"
$ws.Range("F9").Font.Name = "Avenir Next"
$ws.Range("F9").HorizontalAlignment = -4131

$ws.Range("A9").Value = "Πελάτες Τιμή Πώλησης"
$ws.Range("G9").Value = 1.55
$ws.Range("H9").Value = 1.15
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = "Bazaar"
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 2.79

# --- Grand-total row ---------------------------------------------------------
$ws.Range("K35").Value = 352
$ws.Range("L35").Value = 444.04

Write-Output "edits applied"
